# Aggiornamento dati fino al 6 gennaio 2022 (righe 465-491)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(465, 44539, 1, 10, 99.30486593843098),
    @(466, 44540, 2, 11, 109.2353525322741),
    @(467, 44541, 1, 10, 99.30486593843098),
    @(468, 44542, 6, 16, 158.8877855014896),
    @(469, 44543, 8, 19, 188.6792452830189),
    @(470, 44544, 5, 24, 238.3316782522344),
    @(471, 44545, 0, 23, 228.4011916583913),
    @(472, 44546, 9, 31, 307.8450844091361),
    @(473, 44547, 6, 35, 347.5670307845084),
    @(474, 44548, 5, 39, 387.2889771598809),
    @(475, 44550, 9, 42, 417.0804369414101),
    @(476, 44551, 6, 40, 397.2194637537239),
    @(477, 44552, 2, 37, 367.4280039721946),
    @(478, 44553, 6, 43, 427.0109235352532),
    @(479, 44554, 7, 41, 407.149950347567),
    @(480, 44555, 18, 53, 526.3157894736842),
    @(481, 44556, 5, 53, 526.3157894736842),
    @(482, 44557, 11, 55, 546.1767626613704),
    @(483, 44558, 10, 59, 585.8987090367428),
    @(484, 44559, 5, 62, 615.6901688182721),
    @(485, 44560, 4, 60, 595.8291956305859),
    @(486, 44561, 11, 64, 635.5511420059582),
    @(487, 44562, 22, 68, 675.2730883813307),
    @(488, 44563, 7, 70, 695.1340615690168),
    @(489, 44564, 25, 84, 834.1608738828202),
    @(490, 44565, 12, 86, 854.0218470705064),
    @(491, 44566, 23, 104, 1032.770605759682)
)

foreach ($row in $data) {
    $r = $row[0]

    # Replicate the date-column style (s="2": bold, centered, bordered,
    # custom date number format) from the last existing data row (464)
    # onto the new row's column A cell before setting its value.
    $ws.Cells.Item(464, 1).Copy($ws.Cells.Item($r, 1))

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

Write-Host "Done. UsedRange rows:" $ws.UsedRange.Rows.Count
